# Update the cached "datetimeFigureOut" date placeholder text from
# 24.04.2024 to 27.04.2024 across the Slide Master and every Slide
# Layout (these footer-area date fields are stored once per master /
# layout, not per slide).

$p = $ppt.ActivePresentation
$newDate = "27.04.2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePh = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            }
            if ($isDatePh) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder($master.Shapes)

# Every Slide Layout (Custom Layout) under the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder($layout.Shapes)
}
